$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2295.4443
$ws.Range("I62").Value = 1868
$ws.Range("J62").Value = 3266.9092
$ws.Range("K62").Value = 1868
$ws.Range("L62").Value = 3266.9092
$ws.Range("M62").Value = -1244
$ws.Range("N62").Value = -4514.9092

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2295.4443
$ws.Range("I65").Value = 1868
$ws.Range("J65").Value = 3266.9092
$ws.Range("K65").Value = 9340
$ws.Range("L65").Value = 16334.546
$ws.Range("M65").Value = -6220
$ws.Range("N65").Value = -22574.546

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 13278.125
$ws.Range("I132").Value = 13638.044
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 40914.132
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -38384.132
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1798.2727
$ws.Range("I138").Value = 1144.6296
$ws.Range("J138").Value = 2428.5715
$ws.Range("K138").Value = 3433.8888
$ws.Range("L138").Value = 7285.7145
$ws.Range("M138").Value = 1706.1112
$ws.Range("N138").Value = -17565.7145

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 9688.667
$ws.Range("I141").Value = 1099.5555
$ws.Range("J141").Value = 18277.777
$ws.Range("K141").Value = 3298.6665
$ws.Range("L141").Value = 54833.33099999999
$ws.Range("M141").Value = 1881.3335
$ws.Range("N141").Value = -65193.33099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 40811.6
$ws.Range("J33").Value = 48019.332
$ws.Range("L33").Value = 48019.332
$ws.Range("N33").Value = -48677.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 10900
$ws.Range("J107").Value = 10900
$ws.Range("L107").Value = 10900
$ws.Range("N107").Value = -18580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 2000
$ws.Range("K17").Value = 2000
$ws.Range("M17").Value = -1826

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1806.8
$ws.Range("I25").Value = 527.5
$ws.Range("J25").Value = 2126.625
$ws.Range("K25").Value = 527.5
$ws.Range("L25").Value = 2126.625
$ws.Range("M25").Value = -353.5
$ws.Range("N25").Value = -2474.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 40675
$ws.Range("J50").Value = 40675
$ws.Range("L50").Value = 40675
$ws.Range("N50").Value = -41925

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 37072.727
$ws.Range("J51").Value = 37072.727
$ws.Range("L51").Value = 37072.727
$ws.Range("N51").Value = -38544.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 46950
$ws.Range("J55").Value = 46950
$ws.Range("L55").Value = 46950
$ws.Range("N55").Value = -47580

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2148.6785
$ws.Range("I58").Value = 1425.9333
$ws.Range("J58").Value = 2982.6155
$ws.Range("K58").Value = 1425.9333
$ws.Range("L58").Value = 2982.6155
$ws.Range("M58").Value = -1222.9333
$ws.Range("N58").Value = -3388.6155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 42400
$ws.Range("J59").Value = 42400
$ws.Range("L59").Value = 42400
$ws.Range("N59").Value = -44690

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 15693.228
$ws.Range("I60").Value = 20000
$ws.Range("J60").Value = 15488.143
$ws.Range("K60").Value = 20000
$ws.Range("L60").Value = 15488.143
$ws.Range("M60").Value = -19489
$ws.Range("N60").Value = -16510.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 37072.727
$ws.Range("J61").Value = 37072.727
$ws.Range("L61").Value = 37072.727
$ws.Range("N61").Value = -37768.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2491.5454
$ws.Range("I132").Value = 1716.6522
$ws.Range("J132").Value = 4273.8
$ws.Range("K132").Value = 5149.9566
$ws.Range("L132").Value = 12821.4
$ws.Range("M132").Value = -2619.9566
$ws.Range("N132").Value = -17881.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1917.8286
$ws.Range("I134").Value = 1723.3125
$ws.Range("J134").Value = 3992.6667
$ws.Range("K134").Value = 5169.9375
$ws.Range("L134").Value = 11978.0001
$ws.Range("M134").Value = -2634.9375
$ws.Range("N134").Value = -17048.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2148.6785
$ws.Range("I136").Value = 1425.9333
$ws.Range("J136").Value = 2982.6155
$ws.Range("K136").Value = 4277.7999
$ws.Range("L136").Value = 8947.8465
$ws.Range("M136").Value = -1727.7999
$ws.Range("N136").Value = -14047.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1045.96
$ws.Range("I5").Value = 567.6
$ws.Range("J5").Value = 1763.5
$ws.Range("K5").Value = 1702.8
$ws.Range("L5").Value = 5290.5
$ws.Range("M5").Value = -1590.8
$ws.Range("N5").Value = -5514.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2040.2142
$ws.Range("I129").Value = 1005
$ws.Range("J129").Value = 2816.625
$ws.Range("K129").Value = 3015
$ws.Range("L129").Value = 8449.875
$ws.Range("M129").Value = 1985
$ws.Range("N129").Value = -18449.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2307.9102
$ws.Range("I131").Value = 382.14285
$ws.Range("J131").Value = 2729.1719
$ws.Range("K131").Value = 1146.42855
$ws.Range("L131").Value = 8187.5157
$ws.Range("M131").Value = 3893.57145
$ws.Range("N131").Value = -18267.5157

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1045.96
$ws.Range("I135").Value = 567.6
$ws.Range("J135").Value = 1763.5
$ws.Range("K135").Value = 5108.400000000001
$ws.Range("L135").Value = 15871.5
$ws.Range("M135").Value = -2573.400000000001
$ws.Range("N135").Value = -20941.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3423.1924
$ws.Range("I132").Value = 2569.1
$ws.Range("J132").Value = 3957
$ws.Range("K132").Value = 7707.299999999999
$ws.Range("L132").Value = 11871
$ws.Range("M132").Value = -5177.299999999999
$ws.Range("N132").Value = -16931

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7705.2036
$ws.Range("I132").Value = 10617.655
$ws.Range("K132").Value = 31852.965
$ws.Range("M132").Value = -29322.965

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2085.724
$ws.Range("I132").Value = 1767.2632
$ws.Range("J132").Value = 2690.8
$ws.Range("K132").Value = 5301.7896
$ws.Range("L132").Value = 8072.400000000001
$ws.Range("M132").Value = -2771.7896
$ws.Range("N132").Value = -13132.4
